# Update results file with execution times
# Insert a new "Execution Times" worksheet as the first sheet in the workbook
# and populate it with the execution time data, mirroring the target diff.

$wb = $excel.ActiveWorkbook

# Insert the new worksheet before the current first sheet ("5 Nodes") so
# that it becomes, and stays, both the first tab and the active sheet.
$ws = $wb.Sheets.Add($wb.Sheets.Item(1))
$ws.Name = "Execution Times"

# --- Row 1 : column headers (node counts) ------------------------------
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 5
$ws.Range("D1").Value = 10
$ws.Range("E1").Value = "capacity"

# --- Rows 2-3 : "5 nodes" execution times -------------------------------
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 1049
$ws.Range("C2").Value = 503
$ws.Range("D2").Value = 642

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 11299
$ws.Range("C3").Value = 4330
$ws.Range("D3").Value = 3165

# --- Rows 4-5 : "10 nodes" execution times (shaded block) --------------
$ws.Range("A4:D5").Interior.Color = 15921906

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 1934
$ws.Range("C4").Value = 1088
$ws.Range("D4").Value = 760

$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 32206
$ws.Range("C5").Value = 16321
$ws.Range("D5").Value = 4748

# --- Column E labels, merged & centred over each block ------------------
$e2 = $ws.Range("E2")
$e2.Value = "5 nodes"
$e2.VerticalAlignment = -4108
$e2.HorizontalAlignment = -4108
$ws.Range("E2:E3").Merge()

$e4 = $ws.Range("E4")
$e4.Value = "10 nodes"
$e4.VerticalAlignment = -4108
$e4.HorizontalAlignment = -4108
$ws.Range("E4:E5").Merge()

# --- Row 6 / Row 8 labels -----------------------------------------------
$ws.Range("A6").Value = "difficulty"
$ws.Range("A8").Value = "time in seconds"

# Match the selection left on the sheet in the source workbook
$ws.Range("G9").Select()
